# Auto-generated: update computed market/profit columns (H:N) for specific leve rows
# across multiple worksheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 196.66667
$ws.Range("I5").Value = 44.5
$ws.Range("J5").Value = 501
$ws.Range("K5").Value = 44.5
$ws.Range("L5").Value = 501
$ws.Range("M5").Value = 70.5
$ws.Range("N5").Value = -731

$ws.Range("H40").Value = 2361.0454
$ws.Range("I40").Value = 2866.5454
$ws.Range("J40").Value = 1855.5454
$ws.Range("K40").Value = 2866.5454
$ws.Range("L40").Value = 1855.5454
$ws.Range("M40").Value = -2691.5454
$ws.Range("N40").Value = -2205.5454

$ws.Range("H55").Value = 417.27274
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 1000
$ws.Range("N55").Value = -1428

$ws.Range("H113").Value = 2486
$ws.Range("I113").Value = 1619.4546
$ws.Range("J113").Value = 3121.4666
$ws.Range("K113").Value = 1619.4546
$ws.Range("L113").Value = 3121.4666
$ws.Range("M113").Value = 1634.5454
$ws.Range("N113").Value = -9629.4666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1452.8695
$ws.Range("I2").Value = 1356.7222
$ws.Range("J2").Value = 1799
$ws.Range("K2").Value = 1356.7222
$ws.Range("L2").Value = 1799
$ws.Range("M2").Value = -1243.7222
$ws.Range("N2").Value = -2025

$ws.Range("H45").Value = 1908.7941
$ws.Range("I45").Value = 1651.36
$ws.Range("J45").Value = 2623.889
$ws.Range("K45").Value = 1651.36
$ws.Range("L45").Value = 2623.889
$ws.Range("M45").Value = -1274.36
$ws.Range("N45").Value = -3377.889

$ws.Range("H74").Value = 5532.081
$ws.Range("I74").Value = 2548.276
$ws.Range("J74").Value = 16348.375
$ws.Range("K74").Value = 2548.276
$ws.Range("L74").Value = 16348.375
$ws.Range("M74").Value = -1674.276
$ws.Range("N74").Value = -18096.375

$ws.Range("H77").Value = 5532.081
$ws.Range("I77").Value = 2548.276
$ws.Range("J77").Value = 16348.375
$ws.Range("K77").Value = 12741.38
$ws.Range("L77").Value = 81741.875
$ws.Range("M77").Value = -8373.379999999999
$ws.Range("N77").Value = -90477.875

$ws.Range("H116").Value = 1452.8695
$ws.Range("I116").Value = 1356.7222
$ws.Range("J116").Value = 1799
$ws.Range("K116").Value = 1356.7222
$ws.Range("L116").Value = 1799
$ws.Range("M116").Value = 937.2778000000001
$ws.Range("N116").Value = -6387

$ws.Range("H132").Value = 6021.88
$ws.Range("I132").Value = 4866
$ws.Range("J132").Value = 7378.7827
$ws.Range("K132").Value = 14598
$ws.Range("L132").Value = 22136.3481
$ws.Range("M132").Value = -12068
$ws.Range("N132").Value = -27196.3481

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1452.8695
$ws.Range("I3").Value = 1356.7222
$ws.Range("J3").Value = 1799
$ws.Range("K3").Value = 1356.7222
$ws.Range("L3").Value = 1799
$ws.Range("M3").Value = -1242.7222
$ws.Range("N3").Value = -2027

$ws.Range("H134").Value = 3695.7659
$ws.Range("I134").Value = 3766.689
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 11300.067
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -8765.066999999999
$ws.Range("N134").Value = -11370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8134.1665
$ws.Range("I62").Value = 8134.1665
$ws.Range("K62").Value = 8134.1665
$ws.Range("M62").Value = -7510.1665

$ws.Range("H65").Value = 8134.1665
$ws.Range("I65").Value = 8134.1665
$ws.Range("K65").Value = 40670.8325
$ws.Range("M65").Value = -37550.8325

$ws.Range("H132").Value = 2189.4285
$ws.Range("I132").Value = 1772.5264
$ws.Range("J132").Value = 3069.5557
$ws.Range("K132").Value = 5317.5792
$ws.Range("L132").Value = 9208.667099999999
$ws.Range("M132").Value = -2787.5792
$ws.Range("N132").Value = -14268.6671

$ws.Range("H134").Value = 2712.4255
$ws.Range("I134").Value = 1943.875
$ws.Range("K134").Value = 5831.625
$ws.Range("M134").Value = -3296.625

$ws.Range("H135").Value = 42794.668
$ws.Range("J135").Value = 42794.668
$ws.Range("L135").Value = 42794.668
$ws.Range("N135").Value = -52934.668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8922.777
$ws.Range("I80").Value = 12701
$ws.Range("J80").Value = 4200
$ws.Range("K80").Value = 12701
$ws.Range("L80").Value = 4200
$ws.Range("M80").Value = -11703
$ws.Range("N80").Value = -6196

$ws.Range("H83").Value = 8922.777
$ws.Range("I83").Value = 12701
$ws.Range("J83").Value = 4200
$ws.Range("K83").Value = 63505
$ws.Range("L83").Value = 21000
$ws.Range("M83").Value = -58513
$ws.Range("N83").Value = -30984

$ws.Range("H122").Value = 3060.862
$ws.Range("I122").Value = 4058
$ws.Range("J122").Value = 1648.25
$ws.Range("K122").Value = 12174
$ws.Range("L122").Value = 4944.75
$ws.Range("M122").Value = -9724
$ws.Range("N122").Value = -9844.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1058.3334
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1175
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 1175
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -2673

$ws.Range("H71").Value = 1058.3334
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1175
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 5875
$ws.Range("M71").Value = -1256
$ws.Range("N71").Value = -13363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45751.6
$ws.Range("J46").Value = 45751.6
$ws.Range("L46").Value = 45751.6
$ws.Range("N46").Value = -46213.6

$ws.Range("H81").Value = 1560.5714
$ws.Range("I81").Value = 988.6667
$ws.Range("J81").Value = 2590
$ws.Range("K81").Value = 1977.3334
$ws.Range("L81").Value = 5180
$ws.Range("M81").Value = -916.3334
$ws.Range("N81").Value = -7302

$ws.Range("H84").Value = 1560.5714
$ws.Range("I84").Value = 988.6667
$ws.Range("J84").Value = 2590
$ws.Range("K84").Value = 9886.666999999999
$ws.Range("L84").Value = 25900
$ws.Range("M84").Value = -4582.666999999999
$ws.Range("N84").Value = -36508

$ws.Range("H132").Value = 1844.3256
$ws.Range("I132").Value = 1770.48
$ws.Range("J132").Value = 1946.8889
$ws.Range("K132").Value = 5311.440000000001
$ws.Range("L132").Value = 5840.6667
$ws.Range("M132").Value = -2781.440000000001
$ws.Range("N132").Value = -10900.6667

$ws.Range("H134").Value = 45751.6
$ws.Range("J134").Value = 45751.6
$ws.Range("L134").Value = 137254.8
$ws.Range("N134").Value = -142324.8

Write-Output "Updated 168 cells across 27 leve rows."
